# CW3M_McKenzie.xlsx regression-testing workbook update
# Inserts two new model-run rows ("Baseline 2010 C133+" and
# "Baseline 2010-18 C133") into the "2010 and 2010-18" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 26 ("Baseline 2010 C133+"), pushing the
#     existing rows 26-42 down to 27-43. Excel copies the formatting
#     of the row above (25) onto the freshly inserted row, and all
#     relative formulas in the shifted rows are auto-adjusted.
$ws.Rows(26).Insert()

# --- Append a new row at the bottom (44, "Baseline 2010-18 C133").
#     Inserting (rather than just writing past the used range) makes
#     the engine inherit the formatting (incl. style index used by
#     A/B/C/S) from row 43 immediately above it.
$ws.Rows(44).Insert()

# NOTE: new shared-string entries are allocated in the order the
# string values are first written, so B44's label is written before
# B26's to reproduce the original sharedStrings.xml ordering
# (index 60 = "Baseline 2010-18 C133", index 61 = "Baseline 2010 C133+").
$ws.Range("B44").Value2 = "Baseline 2010-18 C133"
$ws.Range("B26").Value2 = "Baseline 2010 C133+"

$ws.Range("A26").Value2 = "CW3M"
$ws.Range("C26").Value2 = 2010
$ws.Range("D26").Value2 = 1090.199341
$ws.Range("E26").Value2 = 1990.4676509999999
$ws.Range("F26").Value2 = 1.2021059999999999
$ws.Range("G26").Value2 = 280.16485599999999
$ws.Range("H26").Value2 = 10.610913999999999
$ws.Range("I26").Value2 = 4.7316050000000001
$ws.Range("J26").Value2 = 8.8404570000000007
$ws.Range("K26").Value2 = 677.36926300000005
$ws.Range("L26").Value2 = 93.229797000000005
$ws.Range("M26").Value2 = 1390.2100829999999
$ws.Range("N26").Value2 = 1207.1236570000001
$ws.Range("O26").Value2 = 6722.8803710000002
$ws.Range("P26").Value2 = 29450.638672000001
$ws.Range("Q26").Value2 = -0.60321499999999995
$ws.Range("R26").Value2 = [double]"-1.7899999999999999E-4"
$ws.Range("S26").Value2 = 2010

$ws.Range("A44").Value2 = "CW3M"
$ws.Range("C44").Value2 = "2010-18"
$ws.Range("D44").Value2 = 1186.9773355555556
$ws.Range("E44").Value2 = 1901.5157334444443
$ws.Range("F44").Value2 = 0.97970299999999988
$ws.Range("G44").Value2 = 280.33542888888883
$ws.Range("H44").Value2 = 9.775355222222224
$ws.Range("I44").Value2 = 5.3870271111111121
$ws.Range("J44").Value2 = 8.145128999999999
$ws.Range("K44").Value2 = 645.93833411111109
$ws.Range("L44").Value2 = 83.47062044444445
$ws.Range("M44").Value2 = 1455.5790065555557
$ws.Range("N44").Value2 = 1191.1918131111113
$ws.Range("O44").Value2 = 4661.9885253333332
$ws.Range("P44").Value2 = 27227.338324888889
$ws.Range("Q44").Value2 = -0.64567966666666665
$ws.Range("R44").Value2 = [double]"-2.0744444444444445E-4"
$ws.Range("S44").Value2 = "2010-18"

# --- Restore the active selection to match the new layout (the sheet
#     view now scrolls a few more rows to keep the same rows visible).
$ws.Range("B27").Select() | Out-Null
